$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at F, shifting the existing "District" column
# (and all its data) from F to G.
$ws.Columns("F").Insert()

# New header for the inserted column.
$ws.Range("F2").Value = "Address"

# Fill the new "Address" column with the school/address portion that was
# embedded in the "Names" (B) column, for every data row except row 20
# (a wrapped continuation line with no address segment of its own).
$ws.Range("F3").Value = "Govt High School Mugalolli"
$ws.Range("F4").Value = "Govt Urdu High School Jamkhandi"
$ws.Range("F5").Value = "AID: Kannada High School Anagawadi,Bilgi"
$ws.Range("F6").Value = "Govt Girls High School"
$ws.Range("F7").Value = "Govt Urdu High School HanagandiJamkhandi"
$ws.Range("F8").Value = "Adarsh Vidyalaya"
$ws.Range("F9").Value = "Govt High SchoolKunchanurJamakandi"
$ws.Range("F10").Value = "Govt High School Muttur(RSMA)Jamakhandi"
$ws.Range("F11").Value = "G H S KulaliMudhol"
$ws.Range("F12").Value = "G H S Simikeri"
$ws.Range("F13").Value = "Govt High School YatnattiBadaradinniBilagi"
$ws.Range("F14").Value = "Govt Urdu High School Bilagi"
$ws.Range("F15").Value = "Shri Basaveshwar High School Hiregulabal"
$ws.Range("F16").Value = "RBG High School LokapurMudhol"
$ws.Range("F17").Value = "Govt High School AdihudiJamakhandi"
$ws.Range("F18").Value = "Govt High School ShirolMudhol"
$ws.Range("F19").Value = "Govt High School Nagur"
$ws.Range("F21").Value = "Govt. Girls High School GuledgudBadami"
$ws.Range("F22").Value = "Govt High School JanamattiBilagi"
$ws.Range("F23").Value = "G H S KudalasangamHunagund"
$ws.Range("F24").Value = "Govt. H P S ADS Tota ShirolMudhol"
$ws.Range("F25").Value = "Govt High SchoolBalakundiHungund"
$ws.Range("F26").Value = "Govt High School TeggiBilagi"
$ws.Range("F27").Value = "Govt High School MannikeriBilagi"
$ws.Range("F28").Value = "Govt Urdu High School MahalingpurMudhol"
$ws.Range("F29").Value = "G H S HaligeriBadami"
$ws.Range("F30").Value = "Shri R T D G P U College NeerabudihalBadami"
$ws.Range("F31").Value = "Adarsha Vidyalaya (RMSA) Mudhol"
$ws.Range("F32").Value = "Govt High School ChikkadapurHungund"
$ws.Range("F33").Value = "GHS YallattiJamkhandi"
$ws.Range("F34").Value = "Govt High School JammanakattiBadami"
$ws.Range("F35").Value = "KLESSCP High SchoolMahaligpurMudhol"
$ws.Range("F36").Value = "G H S Kadampur"
$ws.Range("F37").Value = "Govt High School MuradiHunagund"
$ws.Range("F38").Value = "GHS Alagur (RC)jamakhandi"
$ws.Range("F39").Value = "G H S MadabhaviMudhol"
$ws.Range("F40").Value = "Govt High SchoolKulahalliJamakhandi"
$ws.Range("F41").Value = "Govt High School HiresinganaguttiHunagund"
$ws.Range("F42").Value = "Govt. Adarsha Vidyalaya Jamkhandi"
$ws.Range("F43").Value = "AIDSree Ramalingeshwara High School KundaragiBilagi"
$ws.Range("F44").Value = "Govt High School Kataraki(RSMA)Bilagi"
$ws.Range("F45").Value = "Govt Girls High School Bilagi"
